# aura_光环表: drop the unused "脚本ID/sid" column and add two new auras
# (吸血光环 / Lifesteal Aura, 掉防光环 / Defense-drop Aura) with 5 levels each.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Remove column H (脚本ID / sid) entirely - everything right of it
#    (参数1..参数4, 描述) shifts one column to the left (I->H ... M->L).
$ws.Columns.Item(8).Delete()

# 2) Copy the row-19 formatting pattern down into the 10 new rows (20-29)
#    so the text columns (C name, F icon, G model, L desc) pick up the
#    same "text" style (s="3") used by every other data row.
$ws.Range("A19:L19").Copy()
$ws.Range("A20:L29").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# 3) New aura: 吸血光环 (Lifesteal Aura) - id 2004, levels 1-5
$ws.Range("A20").Value = 2004
$ws.Range("B20").Value = 1
$ws.Range("C20").Value = "吸血光环"
$ws.Range("D20").Value = 4000
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = "none"
$ws.Range("G20").Value = "waitting"
$ws.Range("H20").Value = 10
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = "吸血10%造成的伤害"

$ws.Range("A21").Value = 2004
$ws.Range("B21").Value = 2
$ws.Range("C21").Value = "吸血光环"
$ws.Range("D21").Value = 4000
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = "none"
$ws.Range("G21").Value = "waitting"
$ws.Range("H21").Value = 20
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = "吸血20%造成的伤害"

$ws.Range("A22").Value = 2004
$ws.Range("B22").Value = 3
$ws.Range("C22").Value = "吸血光环"
$ws.Range("D22").Value = 4000
$ws.Range("E22").Value = 0
$ws.Range("F22").Value = "none"
$ws.Range("G22").Value = "waitting"
$ws.Range("H22").Value = 30
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = "吸血30%造成的伤害"

$ws.Range("A23").Value = 2004
$ws.Range("B23").Value = 4
$ws.Range("C23").Value = "吸血光环"
$ws.Range("D23").Value = 4000
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = "none"
$ws.Range("G23").Value = "waitting"
$ws.Range("H23").Value = 40
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = "吸血40%造成的伤害"

$ws.Range("A24").Value = 2004
$ws.Range("B24").Value = 5
$ws.Range("C24").Value = "吸血光环"
$ws.Range("D24").Value = 4000
$ws.Range("E24").Value = 0
$ws.Range("F24").Value = "none"
$ws.Range("G24").Value = "waitting"
$ws.Range("H24").Value = 50
$ws.Range("I24").Value = 0
$ws.Range("J24").Value = 0
$ws.Range("K24").Value = 0
$ws.Range("L24").Value = "吸血50%造成的伤害"

# 4) New aura: 掉防光环 (Defense-drop Aura) - id 2005, levels 1-5
$ws.Range("A25").Value = 2005
$ws.Range("B25").Value = 1
$ws.Range("C25").Value = "掉防光环"
$ws.Range("D25").Value = 3000
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = "none"
$ws.Range("G25").Value = "waitting"
$ws.Range("H25").Value = 7
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = "掉7点防御"

$ws.Range("A26").Value = 2005
$ws.Range("B26").Value = 2
$ws.Range("C26").Value = "掉防光环"
$ws.Range("D26").Value = 3000
$ws.Range("E26").Value = 0
$ws.Range("F26").Value = "none"
$ws.Range("G26").Value = "waitting"
$ws.Range("H26").Value = 9
$ws.Range("I26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = "掉9点防御"

$ws.Range("A27").Value = 2005
$ws.Range("B27").Value = 3
$ws.Range("C27").Value = "掉防光环"
$ws.Range("D27").Value = 3000
$ws.Range("E27").Value = 0
$ws.Range("F27").Value = "none"
$ws.Range("G27").Value = "waitting"
$ws.Range("H27").Value = 11
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = "掉11点防御"

$ws.Range("A28").Value = 2005
$ws.Range("B28").Value = 4
$ws.Range("C28").Value = "掉防光环"
$ws.Range("D28").Value = 3000
$ws.Range("E28").Value = 0
$ws.Range("F28").Value = "none"
$ws.Range("G28").Value = "waitting"
$ws.Range("H28").Value = 13
$ws.Range("I28").Value = 0
$ws.Range("J28").Value = 0
$ws.Range("K28").Value = 0
$ws.Range("L28").Value = "掉13点防御"

$ws.Range("A29").Value = 2005
$ws.Range("B29").Value = 5
$ws.Range("C29").Value = "掉防光环"
$ws.Range("D29").Value = 3000
$ws.Range("E29").Value = 0
$ws.Range("F29").Value = "none"
$ws.Range("G29").Value = "waitting"
$ws.Range("H29").Value = 15
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = "掉15点防御"

# 5) Match the author's final selection/active cell on the sheet.
$ws.Range("O15").Select()
